$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last row (row 19); shifts dimension to A1:E18
$ws.Rows("19").Delete()

# Update data rows 2-18 with bugfixed naive forecaster values
$ws.Range("A2").Value = 39765
$ws.Range("B2").Value = 2008
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 2009
$ws.Range("E2").Value = -0.4513776153963867

$ws.Range("A3").Value = 40130
$ws.Range("B3").Value = 2009
$ws.Range("C3").Value = -1.324983933426882
$ws.Range("D3").Value = 2010
$ws.Range("E3").Value = -0.8803581938132576

$ws.Range("A4").Value = 40494
$ws.Range("B4").Value = 2010
$ws.Range("C4").Value = -0.3900454704678369
$ws.Range("D4").Value = 2011
$ws.Range("E4").Value = -1.213027585730386

$ws.Range("A5").Value = 40862
$ws.Range("B5").Value = 2011
$ws.Range("C5").Value = -0.2995848153489522
$ws.Range("D5").Value = 2012
$ws.Range("E5").Value = -0.3230872999110068

$ws.Range("A6").Value = 41228
$ws.Range("B6").Value = 2012
$ws.Range("C6").Value = -0.2075757021743008
$ws.Range("D6").Value = 2013
$ws.Range("E6").Value = -0.2793004163246238

$ws.Range("A7").Value = 41592
$ws.Range("B7").Value = 2013
$ws.Range("C7").Value = 0.124712275190686
$ws.Range("D7").Value = 2014
$ws.Range("E7").Value = -0.119752617912039

$ws.Range("A8").Value = 41957
$ws.Range("B8").Value = 2014
$ws.Range("C8").Value = -0.255298189276465
$ws.Range("D8").Value = 2015
$ws.Range("E8").Value = -0.05946205208092747

$ws.Range("A9").Value = 42321
$ws.Range("B9").Value = 2015
$ws.Range("C9").Value = 0.07418514192796266
$ws.Range("D9").Value = 2016
$ws.Range("E9").Value = -0.001680662521774678

$ws.Range("A10").Value = 42689
$ws.Range("B10").Value = 2016
$ws.Range("C10").Value = -0.07611406013281474
$ws.Range("D10").Value = 2017
$ws.Range("E10").Value = -0.1247901924724348

$ws.Range("A11").Value = 43053
$ws.Range("B11").Value = 2017
$ws.Range("C11").Value = -0.191300579729714
$ws.Range("D11").Value = 2018
$ws.Range("E11").Value = -0.05219951976568327

$ws.Range("A12").Value = 43418
$ws.Range("B12").Value = 2018
$ws.Range("C12").Value = 0.0970330232288763
$ws.Range("D12").Value = 2019
$ws.Range("E12").Value = -0.1345737582127748

$ws.Range("A13").Value = 43783
$ws.Range("B13").Value = 2019
$ws.Range("C13").Value = -0.7407518902333265
$ws.Range("D13").Value = 2020
$ws.Range("E13").Value = -0.4363737508290888

$ws.Range("A14").Value = 44159
$ws.Range("B14").Value = 2020
$ws.Range("C14").Value = 0.3056679541520335
$ws.Range("D14").Value = 2021
$ws.Range("E14").Value = -0.514812792200714

$ws.Range("A15").Value = 44525
$ws.Range("B15").Value = 2021
$ws.Range("C15").Value = -1.388491535160907
$ws.Range("D15").Value = 2022
$ws.Range("E15").Value = -2.321721165370549

$ws.Range("A16").Value = 44890
$ws.Range("B16").Value = 2022
$ws.Range("C16").Value = -1.678482969789596
$ws.Range("D16").Value = 2023
$ws.Range("E16").Value = -1.107351089172237

$ws.Range("A17").Value = 45254
$ws.Range("B17").Value = 2023
$ws.Range("C17").Value = -0.5999457276250508
$ws.Range("D17").Value = 2024
$ws.Range("E17").Value = -0.4628630633218611

$ws.Range("A18").Value = 45618
$ws.Range("B18").Value = 2024
$ws.Range("C18").Value = -0.05499271238530445
$ws.Range("D18").Value = 2025
$ws.Range("E18").Value = -0.01934819856548309
